# Update existing rows 2-7 and add new rows 8-10 for the "ECs" sending cluster,
# reflecting the refreshed TPM-based NATMI values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rowData = @{
  2 = @("ECs", "Efna5", "Epha1", "ECs", 1, 0.3333333333333333, 0.05800433333333333, 0.174013, 0.02087975181349295, 0.02087975181349295, 3, 1, 2.211928, 6.635783999999999, 0.174938892641363, 0.1749388926413629, 0.1283014090213333, 1.154712681192, 0.003652680660878947, 0.003652680660878947)
  3 = @("ECs", "Efna5", "Epha1", "FAPs", 1, 0.3333333333333333, 0.05800433333333333, 0.174013, 0.02087975181349295, 0.02087975181349295, 3, 1, 3.864911333333334, 11.594734, 0.3056714815357404, 0.3056714815357404, 0.2241816052824445, 2.017634447542, 0.006382344670928952, 0.006382344670928952)
  4 = @("ECs", "Efna5", "Epha1", "MuSCs", 1, 0.3333333333333333, 0.05800433333333333, 0.174013, 0.02087975181349295, 0.02087975181349295, 3, 1, 6.567164333333333, 19.701493, 0.5193896258228967, 0.5193896258228966, 0.3809239890454444, 3.428315901409, 0.01084472648168505, 0.01084472648168505)
  5 = @("FAPs", "Efna5", "Epha1", "ECs", 3, 1, 1.666083666666667, 4.998251, 0.5997381826733804, 0.5997381826733805, 3, 1, 2.211928, 6.635783999999999, 0.174938892641363, 0.1749388926413629, 3.685257112642666, 33.16731401378399, 0.1049175335516246, 0.1049175335516246)
  6 = @("FAPs", "Efna5", "Epha1", "FAPs", 3, 1, 1.666083666666667, 4.998251, 0.5997381826733804, 0.5997381826733805, 3, 1, 3.864911333333334, 11.594734, 0.3056714815357404, 0.3056714815357404, 6.439265645581555, 57.953390810234, 0.1833228588313247, 0.1833228588313247)
  7 = @("FAPs", "Efna5", "Epha1", "MuSCs", 3, 1, 1.666083666666667, 4.998251, 0.5997381826733804, 0.5997381826733805, 3, 1, 6.567164333333333, 19.701493, 0.5193896258228967, 0.5193896258228966, 10.94144523208255, 98.47300708874299, 0.3114977902904311, 0.3114977902904311)
  8 = @("MuSCs", "Efna5", "Epha1", "ECs", 3, 1, 1.053930333333333, 3.161791, 0.3793820655131266, 0.3793820655131266, 3, 1, 2.211928, 6.635783999999999, 0.174938892641363, 0.1749388926413629, 2.331218014349334, 20.980962129144, 0.06636867842885939, 0.06636867842885938)
  9 = @("MuSCs", "Efna5", "Epha1", "FAPs", 3, 1, 1.053930333333333, 3.161791, 0.3793820655131266, 0.3793820655131266, 3, 1, 3.864911333333334, 11.594734, 0.3056714815357404, 0.3056714815357404, 4.073347289843778, 36.660125608594, 0.1159662780334867, 0.1159662780334867)
  10 = @("MuSCs", "Efna5", "Epha1", "MuSCs", 3, 1, 1.053930333333333, 3.161791, 0.3793820655131266, 0.3793820655131266, 3, 1, 6.567164333333333, 19.701493, 0.5193896258228967, 0.5193896258228966, 6.921333694884779, 62.292003253963, 0.1970471090507805, 0.1970471090507805)
}

foreach ($r in $rowData.Keys) {
  $vals = $rowData[$r]
  for ($c = 1; $c -le $vals.Length; $c++) {
    $ws.Cells.Item([int]$r, $c).Value = $vals[$c - 1]
  }
}
